$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.748.89"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3
$ws.Range("D3").Value = "2.472.48"
$ws.Range("E3").Value = "  -0.80%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.87"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.97"
$ws.Range("E6").Value = "  -0.43%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  +0.82%  "

# Row 9
$ws.Range("E9").Value = "  +3.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0878"
$ws.Range("E10").Value = "  +11.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.89"
$ws.Range("E11").Value = "  +0.61%  "

# Row 12
$ws.Range("E12").Value = "  +0.09%  "

# Row 13
$ws.Range("D13").Value = "2.852.67"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14
$ws.Range("E14").Value = "  +0.60%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.80"
$ws.Range("E15").Value = "  -0.88%  "

# Row 16
$ws.Range("D16").Value = "2.489.70"
$ws.Range("E16").Value = "  -1.83%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("E17").Value = "  +3.93%  "

# Row 18
$ws.Range("D18").Value = "41.722.77"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0969"
$ws.Range("E19").Value = "  +4.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.50"
$ws.Range("E20").Value = "  +2.35%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.32"
$ws.Range("E21").Value = "  -0.02%  "

# Row 22
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.47"
$ws.Range("E22").Value = "  +2.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.14"
$ws.Range("E23").Value = "  +1.68%  "

# Row 24
$ws.Range("E24").Value = "  +0.43%  "

# Row 25
$ws.Range("E25").Value = "  +0.62%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.80"
$ws.Range("E27").Value = "  -1.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +1.74%  "

# Row 29
$ws.Range("E29").Value = "  +1.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.36"
$ws.Range("E30").Value = "  -2.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.11"
$ws.Range("E31").Value = "  -0.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.50"
$ws.Range("E32").Value = "  +0.90%  "

# Row 33
$ws.Range("E33").Value = "  +1.25%  "

# Row 34
$ws.Range("E34").Value = "  +0.14%  "

# Row 35
$ws.Range("E35").Value = "  +2.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.49"
$ws.Range("E36").Value = "  -2.87%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -2.00%  "

# Row 38
$ws.Range("E38").Value = "  +0.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("E39").Value = "  -2.18%  "

# Row 40
$ws.Range("E40").Value = "  -1.95%  "

# Row 41
$ws.Range("E41").Value = "  -3.68%  "

# Row 42
$ws.Range("E42").Value = "  +0.01%  "

# Row 43
$ws.Range("D43").Value = "1.969.41"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("E44").Value = "  +0.32%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.97"
$ws.Range("E45").Value = "  -5.17%  "

# Row 46
$ws.Range("E46").Value = "  -1.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.07"
$ws.Range("E47").Value = "  +2.27%  "

# Row 48
$ws.Range("D48").Value = "2.706.64"
$ws.Range("E48").Value = "  -0.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.48"
$ws.Range("E49").Value = "  +0.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.22"
$ws.Range("E50").Value = "  -1.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.87"
$ws.Range("E51").Value = "  +4.43%  "
